$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to new header cells so they match (bold, border, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows: I and J values per row
$data = @{
    2  = @(1, 5)
    3  = @(1, 5)
    4  = @(1, 6)
    5  = @(1, 5)
    6  = @(1, 5)
    7  = @(1, 5)
    8  = @(1, 6)
    9  = @(1, 4)
    10 = @(1, 3)
    11 = @(1, 6)
    12 = @(7, 8)
    13 = @(7, 9)
    14 = @(3, 4)
    15 = @(3, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
